$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Zone master sheet rebuild: columns re-ordered, four audit columns added
# (cr_by, cr_dtimes, upd_by, upd_dtimes, is_deleted, del_dtimes), is_active
# switched from text "TRUE" to a native boolean, and code/hierarchy values
# refreshed to the latest Guinea master data.
# ---------------------------------------------------------------------------

$createdDate = 45079.577595752315

# Header row
$ws.Range("A1").Value = "code"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "hierarchy_level"
$ws.Range("D1").Value = "hierarchy_level_name"
$ws.Range("E1").Value = "hierarchy_path"
$ws.Range("F1").Value = "parent_zone_code"
$ws.Range("G1").Value = "lang_code"
$ws.Range("H1").Value = "is_active"
$ws.Range("I1").Value = "cr_by"
$ws.Range("J1").Value = "cr_dtimes"
$ws.Range("K1").Value = "upd_by"
$ws.Range("L1").Value = "upd_dtimes"
$ws.Range("M1").Value = "is_deleted"
$ws.Range("N1").Value = "del_dtimes"

# Data rows: code, name, hierarchy_level, hierarchy_level_name, hierarchy_path, parent_zone_code
$rows = @(
  @("GN",  "GUINEE",      0, "PAYS",   "GN",   "NULL"),
  @(1,     "CONAKRY",     1, "REGION", "GN/1", "GN"),
  @(4,     "BOKE",        1, "REGION", "GN/4", "GN"),
  @(3,     "KINDIA",      1, "REGION", "GN/3", "GN"),
  @(8,     "MAMOU",       1, "REGION", "GN/8", "GN"),
  @(2,     "LABE",        1, "REGION", "GN/2", "GN"),
  @(5,     "KANKAN",      1, "REGION", "GN/5", "GN"),
  @(6,     "FARANAH",     1, "REGION", "GN/6", "GN"),
  @(7,     "N'ZEREKORE",  1, "REGION", "GN/7", "GN")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = "fra"
    $ws.Cells.Item($r, 8).Value = $true
    $ws.Cells.Item($r, 9).Value = "superadmin"
    $ws.Cells.Item($r, 10).Value = $createdDate
    $ws.Cells.Item($r, 11).Value = "NULL"
    $ws.Cells.Item($r, 12).Value = "NULL"
    $ws.Cells.Item($r, 13).Value = $false
    $ws.Cells.Item($r, 14).Value = "NULL"
}

# cr_dtimes / upd_dtimes-style column uses the date-time number format
$ws.Range("J2:J10").NumberFormat = "mm:ss.0"

# Selection, as left by the author after the edit
[void]$ws.Range("C13").Select()
